# Apply corrected JF code values to row 10 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C10"  = 6.843478434330044
    "I10"  = 2.080435302990134
    "K10"  = 0.0
    "P10"  = 3.919564697009866
    "R10"  = 2.4098376254321625
    "S10"  = 2.080435302990134
    "T10"  = 2.2295331663710933
    "U10"  = 1.8789429385525764
    "V10"  = 0.0
    "W10"  = 1.7336960858251116
    "X10"  = 4.137128678467276
    "Y10"  = 5.888641917474665
    "Z10"  = 2.4271745201551562
    "AA10" = 0.0
    "AI10" = 5.631627857610814
    "AL10" = 2.9264789928727883
    "AN10" = 2.1567179307664386
    "AO10" = 2.4271745201551562
    "AR10" = 4.571628834190022
    "AS10" = 2.7837161109981436
    "AT10" = 1.3677010946922996
    "AU10" = 2.080435302990134
    "AV10" = 0.0
    "AY10" = 4.503829698515615
    "BA10" = 8.481848004303611
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
